# Apply "added pcos sys review test results" edit to the workbook.
$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Sheet1")

# --- Sheet1: new "weighted score" figures in column G (rows 5-9) ---
$ws1.Range("G5").Value = 27.86
$ws1.Range("G6").Formula = "=10*(0.8667)"
$ws1.Range("G7").Formula = "=20*(0.92)"
$ws1.Range("G8").Formula = "=35*0.78"
$ws1.Range("G9").Formula = "=SUM(G5:G8)"

# --- Sheet1: new recall/precision results table (rows 20-21) ---
$ws1.Range("A20").Value = "Recall "
$ws1.Range("B20").Value = "Precision "
$ws1.Range("C20").Value = "Retrieved Paper "
$ws1.Range("D20").Value = "Included "

$ws1.Range("A21").Value = 1
$ws1.Range("C21").Value = 1925
$ws1.Range("D21").Value = 18
$ws1.Range("B21").Formula = "=D21/C21"

# --- Selection / active sheet bookkeeping ---
# Sheet1 becomes the selected/active sheet with H14 highlighted, taking over
# the "tabSelected" flag previously held by "Included articles".
$ws1.Activate() | Out-Null
$ws1.Range("H14").Select() | Out-Null

Write-Output "Applied pcos sys review test results edit."
